$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 3.25
$ws.Range("AJ3").Value = 6
$ws.Range("AK3").Value = 17
$ws.Range("AN3").Value = 10
